$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a plain-text value into a cell without leaving the cell
# reformatted as a Number (Excel auto-detects numeric-looking strings) and
# without changing the cell's style index: force Text format, write the
# value, then drop the format override so the cell style reverts to the
# sheet default (matches the original inline-string cells, which carry no
# style index).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

function Set-Row($row, $b, $c, $d, $e) {
    if ($b -ne $null) { $ws.Range("B$row").Value = $b }
    if ($c -ne $null) { $ws.Range("C$row").Value = $c }
    if ($d -ne $null) { Set-TextValue $ws.Range("D$row") $d }
    if ($e -ne $null) { $ws.Range("E$row").Value = $e }
}

Set-Row 2 $null $null "48.070.60" "  -0.24%  "
Set-Row 3 $null $null "2.488.52" "  -1.45%  "
Set-Row 4 $null $null $null "  -0.06%  "
Set-Row 5 $null $null "317.43" "  -2.02%  "
Set-Row 6 $null $null "105.64" "  -3.10%  "
Set-Row 7 $null $null "0.520" "  -1.67%  "
Set-Row 8 $null $null "0.999" "  +0.02%  "
Set-Row 9 $null $null "0.538" "  -3.21%  "
Set-Row 10 $null $null "38.98" "  -4.78%  "
Set-Row 11 $null $null "20.22" "  -1.27%  "
Set-Row 12 $null $null "0.0801" "  -3.00%  "
Set-Row 13 $null $null $null "  +0.12%  "
Set-Row 15 $null $null "2.880.56" "  -1.46%  "
Set-Row 16 $null $null "2.493.83" "  -1.25%  "
Set-Row 17 $null $null $null "  -3.57%  "
Set-Row 18 $null $null "47.973.44" "  -0.03%  "
Set-Row 19 $null $null "3.00" "  +11.11%  "
Set-Row 20 $null $null $null "  -3.85%  "
Set-Row 21 $null $null "6.57" "  -0.97%  "
Set-Row 22 $null $null "0.0₃0930" "  -2.16%  "
Set-Row 23 $null $null "71.03" "  -1.69%  "
Set-Row 24 $null $null "270.56" "  +0.38%  "
Set-Row 25 $null $null "2.51" "  -2.80%  "
Set-Row 26 $null $null $null "  +0.17%  "
Set-Row 27 $null $null $null "  -1.90%  "
Set-Row 28 $null $null "2.25" "  +1.64%  "
Set-Row 29 $null $null $null "  -4.22%  "
Set-Row 30 $null $null $null "  -3.03%  "
Set-Row 31 $null $null "34.60" "  -3.11%  "
Set-Row 32 $null $null "49.30" "  -0.62%  "
Set-Row 33 $null $null $null "  -0.04%  "
Set-Row 34 $null $null "19.09" "  -4.09%  "
Set-Row 35 $null $null $null "  -2.48%  "
Set-Row 36 $null $null $null "  -2.79%  "
Set-Row 37 $null $null $null "  -2.57%  "
Set-Row 38 $null $null "4.57" "  -3.51%  "
Set-Row 39 $null $null "2.87" "  -4.74%  "
Set-Row 40 $null $null "122.49" "  +2.20%  "
Set-Row 41 $null $null "0.111" "  -1.67%  "
Set-Row 42 "WEMIXToken" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix" "2.22" "  +1.38%  "
Set-Row 43 "EnergySwap" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" "22.15" "  -0.29%  "
Set-Row 44 $null $null $null "  +0.72%  "
Set-Row 45 $null $null "1.998.73" $null
Set-Row 46 $null $null "3.18" "  +0.43%  "
Set-Row 47 $null $null "1.90" "  +0.09%  "
Set-Row 48 $null $null $null "  -2.77%  "
Set-Row 49 $null $null "8.91" "  -2.62%  "
Set-Row 50 $null $null $null "  -1.47%  "
Set-Row 51 $null $null "78.70" "  -1.03%  "
